# Rename the inline picture shapes in the document's headers/footers.
#
# The document has a single section whose footers both show the Pearson
# Edexcel logo (currently named "image1.png") and whose header shows the
# BTEC logo (currently named "image2.jpg"). This commit simply renames
# those inline shapes:
#   - Footer logo(s) (Pearson Edexcel "PearsonLogo.png"): image1.png -> image2.png
#   - Header logo (BTec_Logo-Orange):                     image2.jpg -> image1.jpg
#
# NOTE: InlineShape.Name's getter isn't reliable here, so the shapes are
# identified by their (reliable) AlternativeText instead, and the new
# Name is then assigned unconditionally.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # Footers: Pearson Edexcel logo, image1.png -> image2.png
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }

    # Headers: BTEC logo, image2.jpg -> image1.jpg
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }
}
